$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G holds "K" (formerly "Strike#"); regenerate values per updated calc.
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 2
